# Auto-update draw results: append the 2025-11-03 Pick 4 result as a new
# row (row 48) at the bottom of the results table on the "Results" sheet.
#
# Columns A (date) and C (phase code) look numeric/date-like to Excel's
# automatic type inference, so they are forced to Text ("@") before the
# value is written (otherwise "2025-11-03" would turn into a date serial
# and "251103" would turn into a plain number). The number format is then
# reset back to the default cell style so the new row's formatting matches
# the rest of the table (all other columns never look like numbers/dates,
# so they can be written directly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 48

# A: Date (text, e.g. "2025-11-03") - force text so it isn't parsed as a date
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2025-11-03"
$ws.Range("A$newRow").Style = "Normal"

# B: Game name
$ws.Range("B$newRow").Value = "Pick 4"

# C: Phase code (text, e.g. "251103") - force text so it isn't parsed as a number
$ws.Range("C$newRow").NumberFormat = "@"
$ws.Range("C$newRow").Value = "251103"
$ws.Range("C$newRow").Style = "Normal"

# D: Result
$ws.Range("D$newRow").Value = "6-1-4-5"

# E: InsertedAt timestamp
$ws.Range("E$newRow").Value = "2025-11-03T21:39:02.513+04:00"
